$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) held strings like "5-31-2012-13" which mixed the
# month-day with the season label and were off by a day versus the actual
# NBA game date. Replace them with the corrected ISO-style date
# "2013-05-31" for every data row (BF2:BF31).
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq "5-31-2012-13") {
        $cell.Value2 = "2013-05-31"
    }
}
